$wb = $excel.ActiveWorkbook

# NOTE 1: sheet names "Vector_bf" and "Vector_BF" differ only by case, and
# Worksheets.Item(name) resolves case-insensitively (both names would hit
# the same sheet) - so every sheet below is addressed by its 1-based index
# instead of its name.
#   1 Funciones_Objetivo
#   2 Restricciones_del_lider
#   3 Restricciones_del_follower
#   4 Punto_modificado
#   5 Vector_bf
#   6 Vector_BF
#   7 Vector_Alpha
#
# NOTE 2: several cells hold number-looking text labels (e.g. "0.51",
# "-1", "0") that must stay text, matching the original shared-string
# type, instead of being auto-coerced to real numbers. Briefly flipping
# NumberFormat to "@" (Text) before the assignment forces that, and
# ClearFormats() afterwards removes the quote-prefix style again so no
# stray cell formatting is left behind. (NumberFormat/ClearFormats only
# reliably affect every cell when applied to one contiguous range at a
# time, not a non-contiguous Union, hence the range-by-range calls.)

# --- Sheet 2: Restricciones_del_lider ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A2").Value = "2.3000000000000003 - x"
$ws2.Range("A3").Value = "-2.3000000000000003 + x"

$ws2.Range("B2:B3").NumberFormat = "@"
$ws2.Range("B2").Value = "-3.3000000000000003"
$ws2.Range("B3").Value = "1.3000000000000003"
$ws2.Range("B2:B3").ClearFormats()

$ws2.Range("D2:D3").NumberFormat = "@"
$ws2.Range("D2").Value = "0.51"
$ws2.Range("D3").Value = "0.17"
$ws2.Range("D2:D3").ClearFormats()

# --- Sheet 3: Restricciones_del_follower ---
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("A2").NumberFormat = "@"
$ws3.Range("A2").Value = "0"
$ws3.Range("A2").ClearFormats()

$ws3.Range("D2:F2").NumberFormat = "@"
$ws3.Range("D2").Value = "0.82"
$ws3.Range("E2").Value = "3.4000000000000004"
$ws3.Range("F2").Value = "0"
$ws3.Range("D2:F2").ClearFormats()

$ws3.Range("A3").Value = "-4.761499999999999 + 1.0699999999999998y"

$ws3.Range("B3").NumberFormat = "@"
$ws3.Range("B3").Value = "3.761499999999999"
$ws3.Range("B3").ClearFormats()

$ws3.Range("D3:F3").NumberFormat = "@"
$ws3.Range("D3").Value = "0.81"
$ws3.Range("E3").Value = "6.0"
$ws3.Range("F3").Value = "7.199999999999999"
$ws3.Range("D3:F3").ClearFormats()

# --- Sheet 4: Punto_modificado ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("A2:B2").NumberFormat = "@"
$ws4.Range("A2").Value = "2.3000000000000003"
$ws4.Range("B2").Value = "4.449999999999999"
$ws4.Range("A2:B2").ClearFormats()

# --- Sheet 5: Vector_bf ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("A2").NumberFormat = "@"
$ws5.Range("A2").Value = "1.0653000000000041"
$ws5.Range("A2").ClearFormats()

# --- Sheet 6: Vector_BF ---
$ws6 = $wb.Worksheets.Item(6)
$ws6.Range("A2:A3").NumberFormat = "@"
$ws6.Range("A2").Value = "22.979400000000005"
$ws6.Range("A3").Value = "-10.181"
$ws6.Range("A2:A3").ClearFormats()

# --- Sheet 7: Vector_Alpha ---
$ws7 = $wb.Worksheets.Item(7)
$ws7.Range("A2").Value = 2.07

Write-Output "edit applied"
